{"js": "// Insert two new achievement bullet paragraphs into the\n// \"KEY ACHIEVEMENTS AND IMPACT\" section, right after the existing\n// \"\u2022 Expert methodology validated at highest judicial level\" bullet\n// and right before the \"TECHNICAL SKILLS\" heading.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst anchorText = \"\u2022 Expert methodology validated at highest judicial level\";\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// First new paragraph: plain text bullet.\nconst newPara1 = anchorParagraph.insertParagraph(\n  \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  Word.InsertLocation.after\n);\n\n// Second new paragraph: bullet with a bold, colored \"178%\" run in the middle.\nconst newPara2 = newPara1.insertParagraph(\"\u2022 \", Word.InsertLocation.after);\n\nconst boldRange = newPara2\n  .getRange(Word.RangeLocation.end)\n  .insertText(\"178%\", Word.InsertLocation.end);\nboldRange.font.bold = true;\nboldRange.font.color = \"#2C3E50\";\n\nnewPara2\n  .getRange(Word.RangeLocation.end)\n  .insertText(\" accuracy improvement in racial classification algorithms\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Insert two new achievement bullet paragraphs into the\n# \"KEY ACHIEVEMENTS AND IMPACT\" section, right after the existing\n# \"Expert methodology validated at highest judicial level\" bullet\n# and right before the \"TECHNICAL SKILLS\" heading.\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Expert methodology validated at highest judicial level\")\nif (-not $found) {\n  throw \"Could not find anchor paragraph text\"\n}\n\n# Collapse to the end of the found text (end of that paragraph's content).\n$rng.Collapse(0)\n\n# --- First new paragraph: plain bullet text ---\n$rng.InsertParagraphAfter()\n$rng.Move(4, 1) | Out-Null\n$rng.Text = \"$bullet Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n$rng.Collapse(0)\n\n# --- Second new paragraph: bullet with a bold, colored \"178%\" run ---\n$rng.InsertParagraphAfter()\n$rng.Move(4, 1) | Out-Null\n$rng.Text = \"$bullet \"\n$rng.Collapse(0)\n\n$boldStart = $rng.Start\n$rng.InsertAfter(\"178%\")\n$boldRange = $d.Range($boldStart, $boldStart + 4)\n$boldRange.Font.Bold = 1\n$boldRange.Font.Color = 0x503E2C   # BGR encoding of RGB 2C3E50\n\n$rng.SetRange($boldRange.End, $boldRange.End)\n$rng.InsertAfter(\" accuracy improvement in racial classification algorithms\")\n"}
